# Update res_line / pl_mw values for the 380 kV case (rows 2-25, cols B-L excluding E)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$colB = @(0.7699711380781764, 0.7353886737641915, 0.7144876610159656, 0.7060545776736831, 0.7046593760102269, 0.71437358762023, 0.7579783708798686, 0.8461074605866088, 0.9124339081292874, 0.9429465887147899, 0.9545494412377309, 0.9520484178165134, 0.9439001965027387, 0.9389154563603199, 0.9104466325673286, 0.8930687047293588, 0.8831054599947379, 0.879737604375407, 0.8949152989912079, 0.9462922187826166, 0.9801516271300272, 0.9620546504543483, 0.8940803678430882, 0.8219877070624477)
$colC = @(0.04922058022340536, 0.04617259397027595, 0.04428008751497714, 0.04350357955670603, 0.04337432089766935, 0.04426963671148343, 0.04817399387237487, 0.05566424351540888, 0.06106744206546466, 0.06350414120169035, 0.06442381703116951, 0.06422588391426132, 0.06357986456536935, 0.06318376212298915, 0.06090777104894585, 0.0595060872693125, 0.05869787601189103, 0.05842388592978409, 0.05965550583542267, 0.06376969882033734, 0.06644079099849876, 0.06501680303232149, 0.05958796106908437, 0.05365557693571077)
$colD = @(0.092481901576285, 0.09119600002274808, 0.09043738152759317, 0.09013606622774972, 0.09008650731932022, 0.09043328612244039, 0.09203213262835419, 0.09541091368141963, 0.09803941765790114, 0.09926645781003884, 0.09973556464169064, 0.0996343368294319, 0.0993049625481035, 0.0991037897559508, 0.09795985362380577, 0.09726607176377655, 0.09686997782288387, 0.09673637560692327, 0.09733962097292448, 0.09940158729655479, 0.1007751358419711, 0.1000396907411414, 0.09730636074253596, 0.09447103996460271)
$colF = @(2.425831029891569, 2.429972555645008, 2.433553216909061, 2.435273464905023, 2.435574884958839, 2.433575359399534, 2.427043728704248, 2.422464748366949, 2.424114426428389, 2.425953192681945, 2.426805898277138, 2.426615298225144, 2.426020211239887, 2.425676068997063, 2.424016150430802, 2.423276521683817, 2.422953550166596, 2.422861795194137, 2.423344654600982, 2.426190758682907, 2.428962455141402, 2.427399767919354, 2.423313533190353, 2.422822777307076)
$colG = @(1.729095195964746, 1.731814792600233, 1.734309543445391, 1.735533583150655, 1.735749360376616, 1.734325211515497, 1.729861687113001, 1.72765664344648, 1.73003448075454, 1.731985806652361, 1.732849854085515, 1.732658200031793, 1.732054384309336, 1.731700826485735, 1.729924462309754, 1.729057504954213, 1.728640704058463, 1.728513638127268, 1.729141322360064, 1.732228343272482, 1.734975232953857, 1.733442408185056, 1.729103174255883, 1.727551536448246)
$colH = @(1.532709598043098, 1.539025822164859, 1.543464844494878, 1.545414914723864, 1.5457472495915, 1.543490572246313, 1.534771093126835, 1.522118139958252, 1.515527807586835, 1.513116350340752, 1.512287452691496, 1.512462223883233, 1.513046467966262, 1.513415306516833, 1.51569719808424, 1.517247234735621, 1.5181939848962, 1.518524022819577, 1.517076517170096, 1.5128725748515, 1.510616206281412, 1.511775555826119, 1.517153525337037, 1.525065607205065)
$colI = @(1.558931407012878, 1.565715690066185, 1.57050389180003, 1.572611740783557, 1.572971208435867, 1.570531684780534, 1.561141433105256, 1.547665751336176, 1.540774642113291, 1.538293012354593, 1.537447183166528, 1.537625171170944, 1.538221543463678, 1.538599068147001, 1.540949964625185, 1.542559445325324, 1.543546656404189, 1.543891468888688, 1.542381750477503, 1.538043825888145, 1.53575611489552, 1.536927029453615, 1.5424618934995, 1.550782699081907)
$colJ = @(0.1995054774704883, 0.2002890935999844, 0.2008311523061153, 0.2010673872090791, 0.2011075411331049, 0.200834276103329, 0.1997630396080794, 0.1981446808177125, 0.1972484417691618, 0.1969040289019723, 0.1967826876270919, 0.1968084171176514, 0.1968938642158093, 0.1969473849544769, 0.1972722217571743, 0.1974876941136863, 0.1976175865742107, 0.1976625899345628, 0.1974641401862662, 0.1968685200571443, 0.1965321658619565, 0.1967068492493027, 0.1974747701807438, 0.1985309848058421)
$colK = @(0.5012723984189904, 0.4670836374463363, 0.4462883310801828, 0.4378639671071767, 0.4364681337670504, 0.446174514586545, 0.4894435233016452, 0.5758411030798243, 0.6402491327094992, 0.6697503705125314, 0.6809504086623974, 0.6785370169098144, 0.6706712345526, 0.6658569218909349, 0.6383251781123818, 0.6214866884420189, 0.6118206461274553, 0.6085511691403838, 0.6232772100175339, 0.6729808339315184, 0.7056312883507019, 0.6881900595870718, 0.6224676696112965, 0.552303828945611)
$colL = @(0.364991228517269, 0.3605543007900565, 0.3579862780355398, 0.3569792028204049, 0.3568143628496827, 0.357972536517579, 0.3634290118155121, 0.3753647294920484, 0.3848826899176885, 0.3893743979871545, 0.3910984787737419, 0.3907261385870413, 0.3895157754975287, 0.3887774068648184, 0.3845923947519481, 0.3820664177547854, 0.3806287922166689, 0.3801446596623634, 0.382333734899504, 0.3898706606335907, 0.3949314659971037, 0.392218106443508, 0.3822128353572225, 0.3720039863548692)

$colIndices = @(2, 3, 4, 6, 7, 8, 9, 10, 11, 12)
$data = @($colB, $colC, $colD, $colF, $colG, $colH, $colI, $colJ, $colK, $colL)
$startRow = 2

for ($c = 0; $c -lt $colIndices.Length; $c++) {
    $colValues = $data[$c]
    for ($i = 0; $i -lt $colValues.Length; $i++) {
        $ws.Cells.Item($startRow + $i, $colIndices[$c]).Value = $colValues[$i]
    }
}

